# Update the data values for rows 2-40 (time, CPU Cycles, Input data size)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.968302523621997
$ws.Range("C2").Value = 6.454292139541824
$ws.Range("D2").Value = 0.4494569213310362
$ws.Range("B3").Value = 9.252023311640137
$ws.Range("C3").Value = 6.282854970324147
$ws.Range("D3").Value = 0.3266853406745545
$ws.Range("B4").Value = 15.47387749417902
$ws.Range("C4").Value = 9.814980236431094
$ws.Range("D4").Value = 0.5064586430414989
$ws.Range("B5").Value = 15.76923965739471
$ws.Range("C5").Value = 11.76779120894568
$ws.Range("D5").Value = 0.3604445882999941
$ws.Range("B6").Value = 19.50032959794995
$ws.Range("C6").Value = 5.384490624648259
$ws.Range("D6").Value = 0.1994159906107115
$ws.Range("B7").Value = 23.64428567371602
$ws.Range("C7").Value = 10.72708028679793
$ws.Range("D7").Value = 0.4874454531970882
$ws.Range("B8").Value = 24.79581714028404
$ws.Range("C8").Value = 6.11435447822662
$ws.Range("D8").Value = 0.1097164019370322
$ws.Range("B9").Value = 29.1178153016956
$ws.Range("C9").Value = 1.659738961791062
$ws.Range("D9").Value = 0.149368883682171
$ws.Range("B10").Value = 29.89439917752741
$ws.Range("C10").Value = 1.817738868631651
$ws.Range("D10").Value = 0.1617671317029675
$ws.Range("B11").Value = 30.85529943446785
$ws.Range("C11").Value = 6.885633477717436
$ws.Range("D11").Value = 0.458006571178864
$ws.Range("B12").Value = 33.16401428388629
$ws.Range("C12").Value = 9.730232403613652
$ws.Range("D12").Value = 0.3297553382756093
$ws.Range("B13").Value = 38.76031647558236
$ws.Range("C13").Value = 5.603799047899812
$ws.Range("D13").Value = 0.170441410581006
$ws.Range("B14").Value = 39.55115305736447
$ws.Range("C14").Value = 4.761553316583797
$ws.Range("D14").Value = 0.3709903827184871
$ws.Range("B15").Value = 39.80777509365937
$ws.Range("C15").Value = 2.481549901860521
$ws.Range("D15").Value = 0.2267831297300842
$ws.Range("B16").Value = 42.55025702836069
$ws.Range("C16").Value = 5.717292976747705
$ws.Range("D16").Value = 0.1704305368469851
$ws.Range("B17").Value = 43.59626178713856
$ws.Range("C17").Value = 4.830041604347083
$ws.Range("D17").Value = 0.266305255286413
$ws.Range("B18").Value = 44.94200099836099
$ws.Range("C18").Value = 7.963160917770097
$ws.Range("D18").Value = 0.298947873568302
$ws.Range("B19").Value = 45.04916448941072
$ws.Range("C19").Value = 10.12291259201468
$ws.Range("D19").Value = 0.5347785275628767
$ws.Range("B20").Value = 45.31982797353774
$ws.Range("C20").Value = 9.855805370637492
$ws.Range("D20").Value = 0.5116685137033247
$ws.Range("B21").Value = 48.94323866682335
$ws.Range("C21").Value = 2.836217830172435
$ws.Range("D21").Value = 0.1145085366686621
$ws.Range("B22").Value = 50.67958651947914
$ws.Range("C22").Value = 3.33035982386281
$ws.Range("D22").Value = 0.2544598249290058
$ws.Range("B23").Value = 51.03266627335613
$ws.Range("C23").Value = 7.284504316992748
$ws.Range("D23").Value = 0.3414084458327378
$ws.Range("B24").Value = 52.82347157906169
$ws.Range("C24").Value = 6.138703008276933
$ws.Range("D24").Value = 0.3002420691404664
$ws.Range("B25").Value = 56.43652450032327
$ws.Range("C25").Value = 2.147670259618144
$ws.Range("D25").Value = 0.1981006634022969
$ws.Range("B26").Value = 59.42789561338819
$ws.Range("C26").Value = 3.672904247438154
$ws.Range("D26").Value = 0.142496811888341
$ws.Range("B27").Value = 68.77260594170538
$ws.Range("C27").Value = 4.258527393176168
$ws.Range("D27").Value = 0.2764380447449578
$ws.Range("B28").Value = 70.76210084893634
$ws.Range("C28").Value = 10.26659091204358
$ws.Range("D28").Value = 0.4706126687254231
$ws.Range("B29").Value = 72.26488238377407
$ws.Range("C29").Value = 8.409295892964687
$ws.Range("D29").Value = 0.3180953487505732
$ws.Range("B30").Value = 77.07379352407794
$ws.Range("C30").Value = 7.134848649913564
$ws.Range("D30").Value = 0.4226524432963901
$ws.Range("B31").Value = 77.16439074120886
$ws.Range("C31").Value = 6.906978970784561
$ws.Range("D31").Value = 0.4231461087209499
$ws.Range("B32").Value = 77.49481543064121
$ws.Range("C32").Value = 5.272611471757612
$ws.Range("D32").Value = 0.2038596968604078
$ws.Range("B33").Value = 82.80882213946811
$ws.Range("C33").Value = 6.870244475330873
$ws.Range("D33").Value = 0.4672303418946005
$ws.Range("B34").Value = 84.10681411513478
$ws.Range("C34").Value = 4.797851332576321
$ws.Range("D34").Value = 0.5223145034175196
$ws.Range("B35").Value = 87.24448915556917
$ws.Range("C35").Value = 7.749948912980603
$ws.Range("D35").Value = 0.4518081259672308
$ws.Range("B36").Value = 94.02970406975554
$ws.Range("C36").Value = 1.837574382007818
$ws.Range("D36").Value = 0.3558148478581901
$ws.Range("B37").Value = 94.66941399940016
$ws.Range("C37").Value = 9.023905547582514
$ws.Range("D37").Value = 0.4101601644347047
$ws.Range("B38").Value = 96.83543765651322
$ws.Range("C38").Value = 2.11129436627449
$ws.Range("D38").Value = 0.3211267499499247
$ws.Range("B39").Value = 97.57618033801869
$ws.Range("C39").Value = 4.155580806150788
$ws.Range("D39").Value = 0.1593040422889569
$ws.Range("B40").Value = 97.89032214815171
$ws.Range("C40").Value = 2.057447036649141
$ws.Range("D40").Value = 0.2003235175444442

# The dataset now only needs rows through 40 (previously extended to row 54).
# Remove the now-unused trailing rows 41:54 entirely (shifting dimension to A1:D40).
$ws.Range("A41:D54").EntireRow.Delete()
